# Auto-generated: refresh market-price derived columns (H-N) across multiple sheets
# as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 168.5
$ws.Range("I6").Value = 168.5
$ws.Range("K6").Value = 505.5
$ws.Range("M6").Value = -393.5
# Row 11
$ws.Range("H11").Value = 176.25
$ws.Range("I11").Value = 176.25
$ws.Range("K11").Value = 176.25
$ws.Range("M11").Value = -36.25
# Row 33
$ws.Range("H33").Value = 264.3
$ws.Range("I33").Value = 87.25
$ws.Range("K33").Value = 87.25
$ws.Range("M33").Value = 141.75
# Row 40
$ws.Range("H40").Value = 2108.3333
$ws.Range("J40").Value = 2144.4443
$ws.Range("L40").Value = 2144.4443
$ws.Range("N40").Value = -2494.4443
# Row 64
$ws.Range("H64").Value = 3377.7
$ws.Range("I64").Value = 3557.6
$ws.Range("J64").Value = 3197.8
$ws.Range("K64").Value = 3557.6
$ws.Range("L64").Value = 3197.8
$ws.Range("M64").Value = -3309.6
$ws.Range("N64").Value = -3693.8
# Row 67
$ws.Range("H67").Value = 3377.7
$ws.Range("I67").Value = 3557.6
$ws.Range("J67").Value = 3197.8
$ws.Range("K67").Value = 3557.6
$ws.Range("L67").Value = 3197.8
$ws.Range("M67").Value = -2699.6
$ws.Range("N67").Value = -4913.8
# Row 74
$ws.Range("H74").Value = 5601.25
$ws.Range("J74").Value = 6363
$ws.Range("L74").Value = 6363
$ws.Range("N74").Value = -8235
# Row 77
$ws.Range("H77").Value = 5601.25
$ws.Range("J77").Value = 6363
$ws.Range("L77").Value = 31815
$ws.Range("N77").Value = -41175
# Row 100
$ws.Range("H100").Value = 1540.5
$ws.Range("I100").Value = 1136.9
$ws.Range("J100").Value = 2549.5
$ws.Range("K100").Value = 1136.9
$ws.Range("L100").Value = 2549.5
$ws.Range("M100").Value = -595.9000000000001
$ws.Range("N100").Value = -3631.5
# Row 112
$ws.Range("H112").Value = 2300.1904
$ws.Range("J112").Value = 2365.25
$ws.Range("L112").Value = 7095.75
$ws.Range("N112").Value = -9311.75
# Row 113
$ws.Range("H113").Value = 3485.3333
$ws.Range("I113").Value = 3485.3333
$ws.Range("K113").Value = 3485.3333
$ws.Range("M113").Value = -231.3332999999998
# Row 116
$ws.Range("H116").Value = 6005.8335
$ws.Range("I116").Value = 5037.5
$ws.Range("J116").Value = 7216.25
$ws.Range("K116").Value = 5037.5
$ws.Range("L116").Value = 7216.25
$ws.Range("M116").Value = -1595.5
$ws.Range("N116").Value = -14100.25
# Row 132
$ws.Range("H132").Value = 1270.375
$ws.Range("I132").Value = 1306.3478
$ws.Range("K132").Value = 3919.0434
$ws.Range("M132").Value = -1389.0434
# Row 138
$ws.Range("H138").Value = 4318.884
$ws.Range("J138").Value = 4432
$ws.Range("L138").Value = 13296
$ws.Range("N138").Value = -23576

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 1501
$ws.Range("J97").Value = 1494
$ws.Range("L97").Value = 1494
$ws.Range("N97").Value = -2486
# Row 122
$ws.Range("H122").Value = 670508.25
$ws.Range("I122").Value = 1252953.1
$ws.Range("J122").Value = 4857
$ws.Range("K122").Value = 3758859.3
$ws.Range("L122").Value = 14571
$ws.Range("M122").Value = -3756409.3
$ws.Range("N122").Value = -19471
# Row 124
$ws.Range("H124").Value = 65000
$ws.Range("J124").Value = 65000
$ws.Range("L124").Value = 65000
$ws.Range("N124").Value = -74820

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 4921.125
$ws.Range("I99").Value = 4878
$ws.Range("J99").Value = 4993
$ws.Range("K99").Value = 4878
$ws.Range("L99").Value = 4993
$ws.Range("M99").Value = -3380
$ws.Range("N99").Value = -7989

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 616
$ws.Range("J22").Value = 951
$ws.Range("K22").Value = 616
$ws.Range("L22").Value = 951
$ws.Range("M22").Value = -266
$ws.Range("N22").Value = -1651
# Row 70
$ws.Range("H70").Value = 42500
$ws.Range("J70").Value = 42500
$ws.Range("L70").Value = 42500
$ws.Range("N70").Value = -43130
# Row 73
$ws.Range("H73").Value = 42500
$ws.Range("J73").Value = 42500
$ws.Range("L73").Value = 42500
$ws.Range("N73").Value = -44684
# Row 99
$ws.Range("H99").Value = 9736.138000000001
$ws.Range("I99").Value = 6024.6875
$ws.Range("K99").Value = 6024.6875
$ws.Range("M99").Value = -4526.6875
# Row 126
$ws.Range("H126").Value = 9736.138000000001
$ws.Range("I126").Value = 6024.6875
$ws.Range("K126").Value = 18074.0625
$ws.Range("M126").Value = -15604.0625
# Row 132
$ws.Range("H132").Value = 2621.0454
$ws.Range("I132").Value = 1739.6428
$ws.Range("J132").Value = 4163.5
$ws.Range("K132").Value = 5218.928400000001
$ws.Range("L132").Value = 12490.5
$ws.Range("M132").Value = -2688.928400000001
$ws.Range("N132").Value = -17550.5

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 2513906.2
$ws.Range("J29").Value = 14136.272
$ws.Range("L29").Value = 14136.272
$ws.Range("N29").Value = -14716.272
# Row 36
$ws.Range("H36").Value = 6763.857
$ws.Range("I36").Value = 12000
$ws.Range("K36").Value = 12000
$ws.Range("M36").Value = -11515
# Row 70
$ws.Range("H70").Value = 4434.3335
$ws.Range("I70").Value = 4151.5
$ws.Range("K70").Value = 4151.5
$ws.Range("M70").Value = -3881.5
# Row 73
$ws.Range("H73").Value = 4434.3335
$ws.Range("I73").Value = 4151.5
$ws.Range("K73").Value = 4151.5
$ws.Range("M73").Value = -3215.5
# Row 97
$ws.Range("H97").Value = 1185.9166
$ws.Range("J97").Value = 1422.2
$ws.Range("L97").Value = 1422.2
$ws.Range("N97").Value = -2414.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1796.2222
$ws.Range("I7").Value = 1796.2222
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1796.2222
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1684.2222
$ws.Range("N7").ClearContents()
# Row 20
$ws.Range("H20").Value = 29999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 29999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 29999
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -30451
# Row 46
$ws.Range("H46").Value = 3773.5
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 4326.7144
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 4326.7144
$ws.Range("M46").Value = -2811
$ws.Range("N46").Value = -4702.7144
# Row 55
$ws.Range("H55").Value = 1310.3334
$ws.Range("J55").Value = 979.3333
$ws.Range("L55").Value = 979.3333
$ws.Range("N55").Value = -1325.3333
# Row 100
$ws.Range("H100").Value = 3328
$ws.Range("I100").Value = 2020.4
$ws.Range("K100").Value = 2020.4
$ws.Range("M100").Value = -1479.4
# Row 122
$ws.Range("H122").Value = 10247.167
$ws.Range("I122").Value = 9815.091
$ws.Range("K122").Value = 29445.273
$ws.Range("M122").Value = -26995.273
# Row 126
$ws.Range("H126").Value = 1796.2222
$ws.Range("I126").Value = 1796.2222
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5388.6666
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2918.6666
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 4005.7097
$ws.Range("I132").Value = 3364.261
$ws.Range("K132").Value = 10092.783
$ws.Range("M132").Value = -7562.782999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7682.25
$ws.Range("J62").Value = 7794.4
$ws.Range("L62").Value = 7794.4
$ws.Range("N62").Value = -9042.4
# Row 65
$ws.Range("H65").Value = 7682.25
$ws.Range("J65").Value = 7794.4
$ws.Range("L65").Value = 38972
$ws.Range("N65").Value = -45212
# Row 81
$ws.Range("H81").Value = 20000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 40000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -42122
# Row 84
$ws.Range("H84").Value = 20000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 200000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -210608
# Row 96
$ws.Range("H96").Value = 1024
$ws.Range("I96").Value = 1126.4
$ws.Range("J96").Value = 853.3333
$ws.Range("K96").Value = 1126.4
$ws.Range("L96").Value = 853.3333
$ws.Range("M96").Value = 246.5999999999999
$ws.Range("N96").Value = -3599.3333
# Row 113
$ws.Range("H113").Value = 793.4666999999999
$ws.Range("I113").Value = 678.25
$ws.Range("J113").Value = 925.1429000000001
$ws.Range("K113").Value = 2034.75
$ws.Range("L113").Value = 2775.4287
$ws.Range("M113").Value = 135.25
$ws.Range("N113").Value = -7115.4287
# Row 122
$ws.Range("H122").Value = 4679.154
$ws.Range("J122").Value = 1148.5
$ws.Range("L122").Value = 3445.5
$ws.Range("N122").Value = -8345.5

Write-Output "Updated Sheets via scheduled runner"